# Update recomputed TPM-based NATMI ligand-receptor metrics for Serpinf1-Plxdc2
# (columns E:T, rows 2-17) to reflect the new TPM values used by the pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.802375666666666
$ws.Range("H2").Value = 8.407126999999999
$ws.Range("I2").Value = 0.004883830317835578
$ws.Range("J2").Value = 0.004883830317835578
$ws.Range("M2").Value = 1.275993666666667
$ws.Range("N2").Value = 3.827981
$ws.Range("O2").Value = 0.01190973809858134
$ws.Range("P2").Value = 0.01190973809858134
$ws.Range("Q2").Value = 3.575813602287444
$ws.Range("R2").Value = 32.18232242058699
$ws.Range("S2").Value = 0.00005816514000333297
$ws.Range("T2").Value = 0.00005816514000333297

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.802375666666666
$ws.Range("H3").Value = 8.407126999999999
$ws.Range("I3").Value = 0.004883830317835578
$ws.Range("J3").Value = 0.004883830317835578
$ws.Range("O3").Value = 0.7247630712606724
$ws.Range("P3").Value = 0.7247630712606724
$ws.Range("Q3").Value = 217.6049235673995
$ws.Range("R3").Value = 1958.444312106595
$ws.Range("S3").Value = 0.003539619860670499
$ws.Range("T3").Value = 0.003539619860670499

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.802375666666666
$ws.Range("H4").Value = 8.407126999999999
$ws.Range("I4").Value = 0.004883830317835578
$ws.Range("J4").Value = 0.004883830317835578
$ws.Range("M4").Value = 23.49767666666667
$ws.Range("N4").Value = 70.49303
$ws.Range("O4").Value = 0.2193201912641251
$ws.Range("P4").Value = 0.2193201912641252
$ws.Range("Q4").Value = 65.84931731386777
$ws.Range("R4").Value = 592.64385582481
$ws.Range("S4").Value = 0.001071122599409232
$ws.Range("T4").Value = 0.001071122599409232

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.802375666666666
$ws.Range("H5").Value = 8.407126999999999
$ws.Range("I5").Value = 0.004883830317835578
$ws.Range("J5").Value = 0.004883830317835578
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.714852
$ws.Range("N5").Value = 14.144556
$ws.Range("O5").Value = 0.04400699937662105
$ws.Range("P5").Value = 0.04400699937662105
$ws.Range("Q5").Value = 13.21278651673467
$ws.Range("R5").Value = 118.915078650612
$ws.Range("S5").Value = 0.0002149227177525133
$ws.Range("T5").Value = 0.0002149227177525133

# Row 6
$ws.Range("G6").Value = 543.3469646666667
$ws.Range("I6").Value = 0.9469160079809679
$ws.Range("J6").Value = 0.946916007980968
$ws.Range("M6").Value = 1.275993666666667
$ws.Range("N6").Value = 3.827981
$ws.Range("O6").Value = 0.01190973809858134
$ws.Range("P6").Value = 0.01190973809858134
$ws.Range("Q6").Value = 693.3072857172239
$ws.Range("R6").Value = 6239.765571455015
$ws.Range("S6").Value = 0.01127752165640748
$ws.Range("T6").Value = 0.01127752165640748

# Row 7
$ws.Range("G7").Value = 543.3469646666667
$ws.Range("I7").Value = 0.9469160079809679
$ws.Range("J7").Value = 0.946916007980968
$ws.Range("O7").Value = 0.7247630712606724
$ws.Range("P7").Value = 0.7247630712606724
$ws.Range("S7").Value = 0.6862897541701817
$ws.Range("T7").Value = 0.6862897541701818

# Row 8
$ws.Range("G8").Value = 543.3469646666667
$ws.Range("I8").Value = 0.9469160079809679
$ws.Range("J8").Value = 0.946916007980968
$ws.Range("M8").Value = 23.49767666666667
$ws.Range("N8").Value = 70.49303
$ws.Range("O8").Value = 0.2193201912641251
$ws.Range("P8").Value = 0.2193201912641252
$ws.Range("Q8").Value = 12767.39129355209
$ws.Range("R8").Value = 114906.5216419688
$ws.Range("S8").Value = 0.2076777999814477
$ws.Range("T8").Value = 0.2076777999814478

# Row 9
$ws.Range("G9").Value = 543.3469646666667
$ws.Range("I9").Value = 0.9469160079809679
$ws.Range("J9").Value = 0.946916007980968
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.714852
$ws.Range("N9").Value = 14.144556
$ws.Range("O9").Value = 0.04400699937662105
$ws.Range("P9").Value = 0.04400699937662105
$ws.Range("Q9").Value = 2561.800523052563
$ws.Range("R9").Value = 23056.20470747307
$ws.Range("S9").Value = 0.04167093217293095
$ws.Range("T9").Value = 0.04167093217293096

# Row 10
$ws.Range("G10").Value = 25.919625
$ws.Range("H10").Value = 77.758875
$ws.Range("I10").Value = 0.04517133513098909
$ws.Range("J10").Value = 0.0451713351309891
$ws.Range("M10").Value = 1.275993666666667
$ws.Range("N10").Value = 3.827981
$ws.Range("O10").Value = 0.01190973809858134
$ws.Range("P10").Value = 0.01190973809858134
$ws.Range("Q10").Value = 33.073277342375
$ws.Range("R10").Value = 297.659496081375
$ws.Range("S10").Value = 0.0005379787709733264
$ws.Range("T10").Value = 0.0005379787709733265

# Row 11
$ws.Range("G11").Value = 25.919625
$ws.Range("H11").Value = 77.758875
$ws.Range("I11").Value = 0.04517133513098909
$ws.Range("J11").Value = 0.0451713351309891
$ws.Range("O11").Value = 0.7247630712606724
$ws.Range("P11").Value = 0.7247630712606724
$ws.Range("Q11").Value = 2012.663071589375
$ws.Range("R11").Value = 18113.96764430438
$ws.Range("S11").Value = 0.03273851558248076
$ws.Range("T11").Value = 0.03273851558248077

# Row 12
$ws.Range("G12").Value = 25.919625
$ws.Range("H12").Value = 77.758875
$ws.Range("I12").Value = 0.04517133513098909
$ws.Range("J12").Value = 0.0451713351309891
$ws.Range("M12").Value = 23.49767666666667
$ws.Range("N12").Value = 70.49303
$ws.Range("O12").Value = 0.2193201912641251
$ws.Range("P12").Value = 0.2193201912641252
$ws.Range("Q12").Value = 609.05096757125
$ws.Range("R12").Value = 5481.458708141251
$ws.Range("S12").Value = 0.009906985860584423
$ws.Range("T12").Value = 0.009906985860584427

# Row 13
$ws.Range("G13").Value = 25.919625
$ws.Range("H13").Value = 77.758875
$ws.Range("I13").Value = 0.04517133513098909
$ws.Range("J13").Value = 0.0451713351309891
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.714852
$ws.Range("N13").Value = 14.144556
$ws.Range("O13").Value = 0.04400699937662105
$ws.Range("P13").Value = 0.04400699937662105
$ws.Range("Q13").Value = 122.2071957705
$ws.Range("R13").Value = 1099.8647619345
$ws.Range("S13").Value = 0.001987854916950578
$ws.Range("T13").Value = 0.001987854916950578

# Row 14
$ws.Range("G14").Value = 1.737961666666666
$ws.Range("H14").Value = 5.213884999999999
$ws.Range("I14").Value = 0.003028826570207414
$ws.Range("J14").Value = 0.003028826570207415
$ws.Range("M14").Value = 1.275993666666667
$ws.Range("N14").Value = 3.827981
$ws.Range("O14").Value = 0.01190973809858134
$ws.Range("P14").Value = 0.01190973809858134
$ws.Range("Q14").Value = 2.217628079576111
$ws.Range("R14").Value = 19.958652716185
$ws.Range("S14").Value = 0.00003607253119719468
$ws.Range("T14").Value = 0.00003607253119719469

# Row 15
$ws.Range("G15").Value = 1.737961666666666
$ws.Range("H15").Value = 5.213884999999999
$ws.Range("I15").Value = 0.003028826570207414
$ws.Range("J15").Value = 0.003028826570207415
$ws.Range("O15").Value = 0.7247630712606724
$ws.Range("P15").Value = 0.7247630712606724
$ws.Range("Q15").Value = 134.9530043871361
$ws.Range("R15").Value = 1214.577039484225
$ws.Range("S15").Value = 0.002195181647339454
$ws.Range("T15").Value = 0.002195181647339455

# Row 16
$ws.Range("G16").Value = 1.737961666666666
$ws.Range("H16").Value = 5.213884999999999
$ws.Range("I16").Value = 0.003028826570207414
$ws.Range("J16").Value = 0.003028826570207415
$ws.Range("M16").Value = 23.49767666666667
$ws.Range("N16").Value = 70.49303
$ws.Range("O16").Value = 0.2193201912641251
$ws.Range("P16").Value = 0.2193201912641252
$ws.Range("Q16").Value = 40.83806130239444
$ws.Range("R16").Value = 367.54255172155
$ws.Range("S16").Value = 0.0006642828226837542
$ws.Range("T16").Value = 0.0006642828226837544

# Row 17
$ws.Range("G17").Value = 1.737961666666666
$ws.Range("H17").Value = 5.213884999999999
$ws.Range("I17").Value = 0.003028826570207414
$ws.Range("J17").Value = 0.003028826570207415
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 4.714852
$ws.Range("N17").Value = 14.144556
$ws.Range("O17").Value = 0.04400699937662105
$ws.Range("P17").Value = 0.04400699937662105
$ws.Range("Q17").Value = 8.194232040006666
$ws.Range("R17").Value = 73.74808836006
$ws.Range("S17").Value = 0.000133289568987011
$ws.Range("T17").Value = 0.000133289568987011
